# menyempurnakan fitur upload batch
# Rename the "Dosen Penguji" (examiner) header columns to "Dosen Reviewer"
# on the seminar schedule template sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("D1").Value = "Inisial Dosen Reviewer 1"
$ws.Range("E1").Value = "Inisial Dosen Reviewer 2"

$ws.Range("E6").Select()
